$wb = $excel.ActiveWorkbook

# --- Update values on the "Sources" sheet ---
$sources = $wb.Worksheets.Item("Sources")
$sources.Range("H2").Value = 0
$sources.Range("I2").Value = 1000
$sources.Range("J2").Value = 2000
$sources.Range("K2").Value = 4000

$sources.Range("J3").Value = 2000
$sources.Range("K3").Value = 4000

$sources.Range("H4").Value = 1
$sources.Range("I4").Value = 1000
$sources.Range("J4").Value = 4000
$sources.Range("K4").Value = 8000

# Move the Sources sheet's own selection from J4 to K4
$sources.Range("K4").Select()

# --- Update values on the "Restrictions" sheet ---
$restrictions = $wb.Worksheets.Item("Restrictions")
$restrictions.Range("A2").Value = 400

# Make Restrictions the active sheet and move its selection from B2 to A2
$restrictions.Activate()
$restrictions.Range("A2").Select()
